$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G13:G17 values per the strategy update
$ws.Range("G13").Value = 51.6
$ws.Range("G14").Value = 51.6
$ws.Range("G15").Value = 53.6
$ws.Range("G16").Value = 51.6
$ws.Range("G17").Value = 51.6

# Update the selected cell to G16
$ws.Range("G16").Select()
